$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1627.5454
$ws.Cells.Item(2, 9).Value = 867.1667
$ws.Cells.Item(2, 10).Value = 2540
$ws.Cells.Item(2, 11).Value = 867.1667
$ws.Cells.Item(2, 12).Value = 2540
$ws.Cells.Item(2, 13).Value = -754.1667
$ws.Cells.Item(2, 14).Value = -2766

$ws.Cells.Item(6, 8).Value = 239.6
$ws.Cells.Item(6, 9).Value = 48.666668
$ws.Cells.Item(6, 11).Value = 146.000004
$ws.Cells.Item(6, 13).Value = -34.00000399999999

$ws.Cells.Item(9, 8).Value = 189.5
$ws.Cells.Item(9, 9).Value = 190
$ws.Cells.Item(9, 10).Value = 189.4
$ws.Cells.Item(9, 11).Value = 190
$ws.Cells.Item(9, 12).Value = 189.4
$ws.Cells.Item(9, 13).Value = -21
$ws.Cells.Item(9, 14).Value = -527.4

$ws.Cells.Item(12, 8).Value = 608.2857
$ws.Cells.Item(12, 9).Value = 635.8
$ws.Cells.Item(12, 10).Value = 539.5
$ws.Cells.Item(12, 11).Value = 635.8
$ws.Cells.Item(12, 12).Value = 539.5
$ws.Cells.Item(12, 13).Value = -465.8
$ws.Cells.Item(12, 14).Value = -879.5

$ws.Cells.Item(54, 8).Value = 3416.6667
$ws.Cells.Item(54, 9).Value = 8250
$ws.Cells.Item(54, 11).Value = 8250
$ws.Cells.Item(54, 13).Value = -7764

$ws.Cells.Item(74, 8).Value = 5750
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 5750
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 5750
$ws.Cells.Item(74, 13).Value = ""
$ws.Cells.Item(74, 14).Value = -7622

$ws.Cells.Item(77, 8).Value = 5750
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 5750
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 28750
$ws.Cells.Item(77, 13).Value = ""
$ws.Cells.Item(77, 14).Value = -38110

$ws.Cells.Item(81, 8).Value = 35000
$ws.Cells.Item(81, 10).Value = 35000
$ws.Cells.Item(81, 12).Value = 35000
$ws.Cells.Item(81, 14).Value = -36996

$ws.Cells.Item(84, 8).Value = 35000
$ws.Cells.Item(84, 10).Value = 35000
$ws.Cells.Item(84, 12).Value = 105000
$ws.Cells.Item(84, 14).Value = -114984

$ws.Cells.Item(105, 8).Value = 99999.5
$ws.Cells.Item(105, 10).Value = 99999.5
$ws.Cells.Item(105, 12).Value = 99999.5
$ws.Cells.Item(105, 14).Value = -106987.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2799.3333
$ws.Cells.Item(45, 9).Value = 2799.3333
$ws.Cells.Item(45, 11).Value = 2799.3333
$ws.Cells.Item(45, 13).Value = -2422.3333

$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 14).Value = ""

$ws.Cells.Item(61, 8).Value = 996.25
$ws.Cells.Item(61, 9).Value = 996.25
$ws.Cells.Item(61, 11).Value = 996.25
$ws.Cells.Item(61, 13).Value = -784.25

$ws.Cells.Item(74, 8).Value = 926.875
$ws.Cells.Item(74, 9).Value = 926.875
$ws.Cells.Item(74, 11).Value = 926.875
$ws.Cells.Item(74, 13).Value = -52.875

$ws.Cells.Item(77, 8).Value = 926.875
$ws.Cells.Item(77, 9).Value = 926.875
$ws.Cells.Item(77, 11).Value = 4634.375
$ws.Cells.Item(77, 13).Value = -266.375

$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 14).Value = ""

$ws.Cells.Item(105, 8).Value = 36200
$ws.Cells.Item(105, 10).Value = 36200
$ws.Cells.Item(105, 12).Value = 36200
$ws.Cells.Item(105, 14).Value = -43188

$ws.Cells.Item(136, 8).Value = 996.25
$ws.Cells.Item(136, 9).Value = 996.25
$ws.Cells.Item(136, 11).Value = 2988.75
$ws.Cells.Item(136, 13).Value = -438.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(100, 8).Value = 42772.832
$ws.Cells.Item(100, 10).Value = 42772.832
$ws.Cells.Item(100, 12).Value = 42772.832
$ws.Cells.Item(100, 14).Value = -44936.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1998.5
$ws.Cells.Item(58, 9).Value = 1998.5
$ws.Cells.Item(58, 11).Value = 1998.5
$ws.Cells.Item(58, 13).Value = -1795.5

$ws.Cells.Item(99, 8).Value = 2379.8
$ws.Cells.Item(99, 9).Value = 2725
$ws.Cells.Item(99, 11).Value = 2725
$ws.Cells.Item(99, 13).Value = -1227

$ws.Cells.Item(122, 8).Value = 4308.636
$ws.Cells.Item(122, 10).Value = 1533.3334
$ws.Cells.Item(122, 12).Value = 4600.0002
$ws.Cells.Item(122, 14).Value = -9500.0002

$ws.Cells.Item(126, 8).Value = 2379.8
$ws.Cells.Item(126, 9).Value = 2725
$ws.Cells.Item(126, 11).Value = 8175
$ws.Cells.Item(126, 13).Value = -5705

$ws.Cells.Item(136, 8).Value = 1998.5
$ws.Cells.Item(136, 9).Value = 1998.5
$ws.Cells.Item(136, 11).Value = 5995.5
$ws.Cells.Item(136, 13).Value = -3445.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 38402.168
$ws.Cells.Item(99, 9).Value = 47603.25
$ws.Cells.Item(99, 10).Value = 20000
$ws.Cells.Item(99, 11).Value = 47603.25
$ws.Cells.Item(99, 12).Value = 20000
$ws.Cells.Item(99, 13).Value = -45357.25
$ws.Cells.Item(99, 14).Value = -24492

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 3007.25
$ws.Cells.Item(19, 9).Value = 2008.3334
$ws.Cells.Item(19, 11).Value = 2008.3334
$ws.Cells.Item(19, 13).Value = -1838.3334

$ws.Cells.Item(22, 8).Value = 763.4286
$ws.Cells.Item(22, 9).Value = 649.3333
$ws.Cells.Item(22, 10).Value = 849
$ws.Cells.Item(22, 11).Value = 649.3333
$ws.Cells.Item(22, 12).Value = 849
$ws.Cells.Item(22, 13).Value = -354.3333
$ws.Cells.Item(22, 14).Value = -1439

$ws.Cells.Item(27, 8).Value = 763.4286
$ws.Cells.Item(27, 9).Value = 649.3333
$ws.Cells.Item(27, 10).Value = 849
$ws.Cells.Item(27, 11).Value = 649.3333
$ws.Cells.Item(27, 12).Value = 849
$ws.Cells.Item(27, 13).Value = -542.3333
$ws.Cells.Item(27, 14).Value = -1063

$ws.Cells.Item(68, 8).Value = 30000
$ws.Cells.Item(68, 10).Value = 30000
$ws.Cells.Item(68, 12).Value = 30000
$ws.Cells.Item(68, 14).Value = -31498

$ws.Cells.Item(71, 8).Value = 30000
$ws.Cells.Item(71, 10).Value = 30000
$ws.Cells.Item(71, 12).Value = 150000
$ws.Cells.Item(71, 14).Value = -157488

$ws.Cells.Item(95, 8).Value = 16375
$ws.Cells.Item(95, 10).Value = 16375
$ws.Cells.Item(95, 12).Value = 16375
$ws.Cells.Item(95, 14).Value = -21867

$ws.Cells.Item(97, 8).Value = 7000.25
$ws.Cells.Item(97, 10).Value = 7000.25
$ws.Cells.Item(97, 12).Value = 7000.25
$ws.Cells.Item(97, 14).Value = -8982.25

$ws.Cells.Item(101, 8).Value = 26560.166
$ws.Cells.Item(101, 10).Value = 26560.166
$ws.Cells.Item(101, 12).Value = 26560.166
$ws.Cells.Item(101, 14).Value = -33050.166

$ws.Cells.Item(106, 8).Value = 11226.5
$ws.Cells.Item(106, 10).Value = 11226.5
$ws.Cells.Item(106, 12).Value = 11226.5
$ws.Cells.Item(106, 14).Value = -13750.5

$ws.Cells.Item(132, 8).Value = 21665.445
$ws.Cells.Item(132, 9).Value = 21712.572
$ws.Cells.Item(132, 11).Value = 65137.716
$ws.Cells.Item(132, 13).Value = -62607.716

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(94, 8).Value = 26597.6
$ws.Cells.Item(94, 10).Value = 26597.6
$ws.Cells.Item(94, 12).Value = 26597.6
$ws.Cells.Item(94, 14).Value = -28399.6

$ws.Cells.Item(95, 8).Value = 4672
$ws.Cells.Item(95, 10).Value = 4672
$ws.Cells.Item(95, 12).Value = 4672
$ws.Cells.Item(95, 14).Value = -10164

$ws.Cells.Item(97, 8).Value = 54997
$ws.Cells.Item(97, 10).Value = 54997
$ws.Cells.Item(97, 12).Value = 54997
$ws.Cells.Item(97, 14).Value = -56979

$ws.Cells.Item(101, 8).Value = 15040.4
$ws.Cells.Item(101, 10).Value = 15040.4
$ws.Cells.Item(101, 12).Value = 15040.4
$ws.Cells.Item(101, 14).Value = -21530.4

$ws.Cells.Item(105, 8).Value = 17586.25
$ws.Cells.Item(105, 10).Value = 17586.25
$ws.Cells.Item(105, 12).Value = 17586.25
$ws.Cells.Item(105, 14).Value = -24574.25

$ws.Cells.Item(132, 8).Value = 1019.8
$ws.Cells.Item(132, 9).Value = 901.3333
$ws.Cells.Item(132, 10).Value = 1197.5
$ws.Cells.Item(132, 11).Value = 2703.9999
$ws.Cells.Item(132, 12).Value = 3592.5
$ws.Cells.Item(132, 13).Value = -173.9998999999998
$ws.Cells.Item(132, 14).Value = -8652.5
